$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: task renamed (domain files split into Player/Card) and is now Done
$ws.Range("B6").Value = "Create Domain files Player and Card"
$ws.Range("D6").Value = "Done"

# Row 7: existing "Get Api call working" task, now marked Done
$ws.Range("B7").Value = "Get Api call working"
$ws.Range("D7").Value = "Done"

# Row 8: existing "Fix CORS" task, now marked Done
$ws.Range("B8").Value = "Fix CORS"
$ws.Range("D8").Value = "Done"

# Row 9: new task, no status yet
$ws.Range("B9").Value = "Create domain file Resources"

# Row 10: new task, marked Done
$ws.Range("B10").Value = "Learn how to make request that keeps listening"
$ws.Range("D10").Value = "Done"

# Row 11: new task, no status yet
$ws.Range("B11").Value = "Create domain gamestate object"

# Leave selection on the next empty row, matching Excel's behavior after data entry
$ws.Range("B12").Select()
